# Updates the cryptocurrency price/volume table (columns D and E, rows 2-51)
# to reflect the latest scraped values, mirroring the GitHub Actions data
# refresh described in the commit message.
# A leading apostrophe is used for values that look like plain numbers
# (e.g. "374.44") so Excel stores/keeps them as text instead of auto-
# converting them to numeric literals, matching the original cell typing.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '50.992.62'
$ws.Cells.Item(2, 5).Value = '  +0.03%  '
$ws.Cells.Item(3, 4).Value = '2.937.66'
$ws.Cells.Item(3, 5).Value = '  +0.01%  '
$ws.Cells.Item(4, 5).Value = '  -0.06%  '
$ws.Cells.Item(5, 4).Value = '''374.44'
$ws.Cells.Item(5, 5).Value = '  -0.50%  '
$ws.Cells.Item(6, 4).Value = '''101.82'
$ws.Cells.Item(6, 5).Value = '  -2.01%  '
$ws.Cells.Item(7, 5).Value = '  -0.88%  '
$ws.Cells.Item(8, 5).Value = '  +0.11%  '
$ws.Cells.Item(9, 5).Value = '  -0.98%  '
$ws.Cells.Item(10, 4).Value = '''36.46'
$ws.Cells.Item(10, 5).Value = '  -0.54%  '
$ws.Cells.Item(11, 5).Value = '  -0.26%  '
$ws.Cells.Item(12, 4).Value = '''0.0836'
$ws.Cells.Item(12, 5).Value = '  +0.07%  '
$ws.Cells.Item(13, 4).Value = '3.399.72'
$ws.Cells.Item(13, 5).Value = '  -0.28%  '
$ws.Cells.Item(14, 4).Value = '''17.90'
$ws.Cells.Item(14, 5).Value = '  -1.83%  '
$ws.Cells.Item(15, 5).Value = '  -1.07%  '
$ws.Cells.Item(16, 4).Value = '2.944.32'
$ws.Cells.Item(16, 5).Value = '  -0.08%  '
$ws.Cells.Item(17, 4).Value = '''0.976'
$ws.Cells.Item(17, 5).Value = '  +2.88%  '
$ws.Cells.Item(18, 4).Value = '50.911.34'
$ws.Cells.Item(18, 5).Value = '  -0.16%  '
$ws.Cells.Item(19, 4).Value = '''3.15'
$ws.Cells.Item(19, 5).Value = '  -4.87%  '
$ws.Cells.Item(20, 5).Value = '  -1.81%  '
$ws.Cells.Item(21, 4).Value = '''12.58'
$ws.Cells.Item(21, 5).Value = '  -1.51%  '
$ws.Cells.Item(22, 5).Value = '  +0.35%  '
$ws.Cells.Item(23, 4).Value = '''264.00'
$ws.Cells.Item(23, 5).Value = '  +1.60%  '
$ws.Cells.Item(24, 4).Value = '''68.28'
$ws.Cells.Item(24, 5).Value = '  -0.70%  '
$ws.Cells.Item(25, 5).Value = '  +3.50%  '
$ws.Cells.Item(26, 4).Value = '''8.33'
$ws.Cells.Item(26, 5).Value = '  +13.29%  '
$ws.Cells.Item(27, 4).Value = '''7.80'
$ws.Cells.Item(27, 5).Value = '  +9.57%  '
$ws.Cells.Item(28, 5).Value = '  -0.11%  '
$ws.Cells.Item(29, 5).Value = '  -0.03%  '
$ws.Cells.Item(30, 5).Value = '  +1.01%  '
$ws.Cells.Item(31, 5).Value = '  -0.28%  '
$ws.Cells.Item(32, 4).Value = '''9.85'
$ws.Cells.Item(32, 5).Value = '  +0.93%  '
$ws.Cells.Item(33, 4).Value = '''50.82'
$ws.Cells.Item(33, 5).Value = '  -0.05%  '
$ws.Cells.Item(34, 4).Value = '''33.53'
$ws.Cells.Item(34, 5).Value = '  -2.28%  '
$ws.Cells.Item(35, 4).Value = '''0.0448'
$ws.Cells.Item(35, 5).Value = '  +1.30%  '
$ws.Cells.Item(36, 5).Value = '  -2.76%  '
$ws.Cells.Item(37, 5).Value = '  -0.25%  '
$ws.Cells.Item(38, 5).Value = '  -1.83%  '
$ws.Cells.Item(39, 4).Value = '''2.56'
$ws.Cells.Item(39, 5).Value = '  -0.23%  '
$ws.Cells.Item(40, 5).Value = '  -0.24%  '
$ws.Cells.Item(41, 4).Value = '''16.39'
$ws.Cells.Item(41, 5).Value = '  -3.74%  '
$ws.Cells.Item(42, 5).Value = '  -1.96%  '
$ws.Cells.Item(43, 4).Value = '''120.17'
$ws.Cells.Item(43, 5).Value = '  -1.48%  '
$ws.Cells.Item(44, 4).Value = '''0.290'
$ws.Cells.Item(44, 5).Value = '  +4.50%  '
$ws.Cells.Item(45, 4).Value = '''20.96'
$ws.Cells.Item(45, 5).Value = '  -3.77%  '
$ws.Cells.Item(46, 5).Value = '  -1.50%  '
$ws.Cells.Item(47, 4).Value = '''3.24'
$ws.Cells.Item(47, 5).Value = '  +1.90%  '
$ws.Cells.Item(48, 5).Value = '  -3.08%  '
$ws.Cells.Item(49, 4).Value = '1.974.68'
$ws.Cells.Item(49, 5).Value = '  -2.58%  '
$ws.Cells.Item(50, 4).Value = '''0.0341'
$ws.Cells.Item(50, 5).Value = '  -0.19%  '
$ws.Cells.Item(51, 5).Value = '  -0.44%  '
